# kliceni_data.xlsx - "pokus with commit" edit
# Adds a "total" (T) column with the per-person quota (25) and a
# percentage-used (U) column on the "proklicovani" sheet, plus a
# percentage column (K/L) on the little summary table in rows 19-22,
# an incidental format touch on M18, and tidies up the page setup /
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet          # "proklicovani" is already the active sheet

# --- header for the new "total" column -------------------------------
$ws.Range("T1").Value = "total"

# --- per-row quota (T) and percentage used (U), rows 2-13 ------------
$ws.Range("T2").Value = 25
$ws.Range("U2").Formula = "=R2/T2"

$ws.Range("T3:T13").Value = 25
$ws.Range("U3:U13").Formula = "=R3/T3"

# --- stray formatting touch on M18 (General number format) -----------
$ws.Range("M18").NumberFormat = "General"

# --- summary table (rows 19-22): add quota (K) and percentage (L) ----
$ws.Range("K19").Value = 75
$ws.Range("L19").Formula = "=I19/K19"

$ws.Range("K20:K22").Value = 75
$ws.Range("L20:L22").Formula = "=I20/K20"

# --- page setup ---------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- selection moved to N22 -------------------------------------------
$ws.Range("N22").Select() | Out-Null
